$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "add id`tfresh egg rigatoni (shorter cut) 2/5lbs`tcase`t5"
$ws.Range("B16").Value = "5 32RIG1"
$ws.Range("C16").Value = ""
